# Form101 workbook update
# - Insert two new columns (maritalStatus, sex) after HomeZip / before HomePhone
# - Append ten new trailing columns with additional claim-intake fields
# - Add a mailto hyperlink on the new SubmitterEmail cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two columns at J:K. This shifts the existing J..AB columns to
#    L..AD and carries their formatting/values along automatically.
# ---------------------------------------------------------------------------
$ws.Columns("J:K").Insert()
$ws.Columns("J:K").ColumnWidth = 7.42
$ws.Columns("AK:AK").ColumnWidth = 9.67

# ---------------------------------------------------------------------------
# 2. Populate all of the new text cells (new columns J:K plus the appended
#    AE:AN columns). Order matches the order these values were first entered
#    so the shared-string table comes out in the same sequence.
# ---------------------------------------------------------------------------
$ws.Range("AE1").Value = "PartofBodyInjury"
$ws.Range("K1").Value = "sex"
$ws.Range("K2").Value = "Female"
$ws.Range("J1").Value = "maritalStatus"
$ws.Range("J2").Value = "Single"
$ws.Range("AF1").Value = "NatureofInjury"
$ws.Range("AF2").Value = "Specific Injury - Laceration"
$ws.Range("AH1").Value = "WhatHappened"
$ws.Range("AI1").Value = "WhatObject"
$ws.Range("AI2").Value = "Test Object"
$ws.Range("AL1").Value = "Wasworkedemp"
$ws.Range("AM1").Value = "GrossEarnings"
$ws.Range("AN1").Value = "SubmitterEmail"
$ws.Range("AN2").Value = "test@gmail.com"
$ws.Range("AJ1").Value = "Doing"
$ws.Range("AJ2").Value = "TestDoing"
$ws.Range("AK1").Value = "DOLastHire"
$ws.Range("AG1").Value = "CauseofInjury"
$ws.Range("AG2").Value = "Motor Vehicle - Vehicle Upset"
$ws.Range("AE2").Value = "Head - Eyes"

# Cells that reuse already-existing shared strings.
$ws.Range("AH2").Value = "NA"
$ws.Range("AL2").Value = "Yes"

# ---------------------------------------------------------------------------
# 3. AK2 is a date field - copy the date format from an existing date cell
#    (T2) so it reuses the same "short date" style instead of minting a new
#    numFmt, then write the numeric serial directly (Value2-style literal).
# ---------------------------------------------------------------------------
$ws.Range("T2").Copy($ws.Range("AK2"))
$ws.Range("AK2").Value = 40159

# Plain numeric (non-string) cell.
$ws.Range("AM2").Value = 2000

# ---------------------------------------------------------------------------
# 4. Mailto hyperlink on the SubmitterEmail value cell.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("AN2"), "mailto:test@gmail.com") | Out-Null

# ---------------------------------------------------------------------------
# 5. Selection matches the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("AC6").Select()
